# Applies the row 7/8/9 data update described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 7 ----
$ws.Range("A7").Value = 111638335
$ws.Range("B7").Value = 95532
$ws.Range("E7").Value = 221945
$ws.Range("F7").Value = "Revlummer"
$ws.Range("G7").Value = "Lycopodium annotinum"
$ws.Range("H7").Value = "L."
$ws.Range("I7").Value = "'"
$ws.Range("Q7").Value = 595260
$ws.Range("R7").Value = 7035583
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()
$ws.Range("AC7").Value = "Stort bestånd."

# ---- Row 8 ----
$ws.Range("A8").Value = 111638343
$ws.Range("I8").Value = "'15"
$ws.Range("Q8").Value = 595350
$ws.Range("R8").Value = 7035563
$ws.Range("Z8").ClearContents()
$ws.Range("AB8").ClearContents()
$ws.Range("AC8").Value = "I örtrik brant."

# ---- Row 9 ----
$ws.Range("A9").Value = 111638342
$ws.Range("B9").Value = 95723
$ws.Range("E9").Value = 220250
$ws.Range("F9").Value = "Strutbräken"
$ws.Range("G9").Value = "Matteuccia struthiopteris"
$ws.Range("H9").Value = "(L.) Tod."
$ws.Range("I9").Value = "'20"
$ws.Range("Q9").Value = 595216
$ws.Range("R9").Value = 7035642
$ws.Range("Z9").ClearContents()
$ws.Range("AB9").ClearContents()
$ws.Range("AC9").Value = "Långt ner i branten, nedanför lövskogen."
